$d = $word.ActiveDocument

# --- 1. "Ativação: 01/01/2012" -> "Ativação: Semestral" ----------------
$d.Content.Find.Execute("Ativação: 01/01/2012", $false, $false, $false, $false, $false, $true, 1, $false, "Ativação: Semestral", 2) | Out-Null

# --- 2. Objetivos (PT) paragraph -> "01/01/2025" -----------------------
$p6 = $d.Paragraphs(6)
$p6.Range.Text = "01/01/2025"

# --- 3. Objectives (EN, italic) paragraph: merge + append new sentence -
$p7 = $d.Paragraphs(7)
$p7.Range.Text = "Enable the students with the knowledge needed to use the concepts of psychology in the workplace.Promote the recognition of workers' needs, both in its area of operation as in human relationships that have in the company.Promote strategies for the development of transversal skills necessary for the future engineer: effective communication, teamwork, commitment, initiative, responsibility, ethics, among others."

# --- 4. ListBullet paragraph (8188658...) gains a run before + several after
$p9 = $d.Paragraphs(9)
$p9.Range.InsertBefore("Possibilitar aos alunos da disciplina o conhecimento necessário para a utilização dos conceitos da psicologia em ambiente de trabalho.Favorecer o reconhecimento das necessidades dos trabalhadores tanto na sua área de atuação quanto nos relacionamentos humanos que terá na empresa.Promover estratégias para o desenvolvimento de competências transversais necessárias ao futuro engenheiro: comunicação eficaz, trabalho em equipe, comprometimento, iniciativa, responsabilidade, ética, entre outras." + [char]11)
$p9.Range.InsertAfter([char]11 + "1 - Introdução a Psicologia aplicada ao trabalho.  2 - Conceito de Comunicação. 3 - Relações Humanas no Trabalho. 4 - Psicologia nas Organizações de Trabalho: grupos e equipes. 5 - Recrutamento e Seleção. 6 - Motivação. 7 – Liderança.   8 – Treinamento e Desenvolvimento. 9- Avaliação de desempenho." + [char]11 + "1.Introdução: conceituar psicologia como ciência e como aplicação; a psicologia aplicada ao trabalho. A psicologia nas relações humanas no trabalho.2.Conceito de Comunicação: sistemas, funções, axiomas da comunicação humana. Processos de comunicação e o convívio sócio-comunicacional na empresa.3.Relações Humanas no Trabalho: relações humanas em grupos; como participar de um grupo de trabalho, trabalho em equipe, dinâmicas grupais.4.Psicologia nas Organizações de Trabalho: conceitos de organização e de trabalho. Organização e trabalho e sua importância na saúde mental e produtividade do trabalhador: estresse, síndrome de burnout, síndrome de Karoshi; L.E.R.; qualidade de vida; assédios sexual e moral no ambiente de trabalho.5.Recrutamento e Seleção: recrutamento e seleção de pessoal; colocação e acompanhamento; avaliação de desempenho; medidas de avaliação e sua importância na seleção; experiências práticas em sala de aula como facilitadoras do processo de seleção.6.Motivação: as necessidades básicas e psicológicas do ser humano; motivação e conflitos; fatores esquecidos como motivadores na empresa: inveja, ciúme, medo, abuso de poder. Avaliação de motivação.7. Liderança: definição, teorias e desenvolvimento de lideranças8. Treinamento e Desenvolvimento: definição, diferenciação, etapas, dificuldades9. Avaliação de desempenho: definição, tipos, periodicidade, importância" + [char]11 + "Serão aplicadas provas dissertativas com estudo de caso e situações, para levar os alunos à maior reflexão sobre a utilização dos conceitos aprendidos para o futuro engenheiro em seu trabalho cotidiano em empresas.Será solicitada a realização de atividades variadas (avaliação processual) sobre cada tema, com ênfase no desenvolvimento das habilidades transversais. Tais atividades poderão ser: apresentações, elaboração de folder, mapas conceituais e pitch para processo seletivo, relatório e leitura ativa, entrevista com trabalhadores, etc.Será realizada uma atividade extensionista de confecção de currículo para a comunidade. Para tanto, os alunos deverão participar de algum evento que ocorra durante o semestre (feiras de ciências, feira de profissões, reunião com calouros de outras instituições). Caberá ao aluno o planejamento e a execução da atividade, com supervisão da profa." + [char]11 + "P1 = soma de todas as atividades até 8,0 + relatório do desenvolvimento da atividade de extensão.P2 = prova dissertativa aplicada no final do semestre.Média = (P1+ P2)/ 2.")

# --- 5. "Programa resumido" PT paragraph -> "Nova avaliação, similar à P2."
$p11 = $d.Paragraphs(11)
$p11.Range.Text = "Nova avaliação, similar à P2."

# --- 6. "Programa resumido" EN (italic) paragraph: merge + update text -
$p12 = $d.Paragraphs(12)
$p12.Range.Text = "Introduction to Psychology applied to work. Concept of Communication. Human Relations at Work. Psychology in Work Organizations: groups and teams. Recruitment and Selection. Motivation. Leadership. Training and Development. Performance evaluation."

# --- 7. "Programa" PT paragraph -> bibliography text --------------------
$p14 = $d.Paragraphs(14)
$p14.Range.Text = "1.ROBBINS, S. P. Fundamentos do Comportamento Organizacional. 8ª. Ed., São Paulo: Pearson Prentice Hall, 2009.2.CHIAVENATO, I. Gestão de Pessoas: o novo papel da Gestão do Talento Humano. 5ª.ed., São Paulo: Atlas, 2020.3.SHILD, M (trad) A arte de dar feedback. Harvard Business Review. Rio de Janeiro: Sextante, 2019.4.BROWNIE, S. HR on Purpose: Developing Deliberate People Passion — Steve Browne. Alexandria: Society for Human Resource Management, 2017.5.BERGAMINI, C.W. Motivação nas Organizações.7ª.Ed. São Paulo: Atlas, 2018.6.ROSENBERG, M.B. Comunicação Não-Violenta. Rio de Janeiro: Sextante, 2018."

# --- 8. "Programa" EN (italic) paragraph: merge + append ----------------
$p15 = $d.Paragraphs(15)
$p15.Range.Text = "Introduction : conceptualize psychology as science and application; psychology applied to work. The psychology of human relations at work. Concept of Communication : Systems, functions , axioms of human communication. Communication processes and the social and communicative interaction in the company. Human relations at work: the role of masks in human interaction; human relations in groups; how to be a part of a workgroup, teamwork and group dynamics Psychology in Work Organizations : Organization concepts and work. Organization and work and their importance in mental health and worker productivity: stress, burnout , Karoshi syndrome ; L.E.R .; quality of life; sexual and moral harassment in the workplace ; alcohol and drugs at work; mental disorders in the company. Recruitment and Selection: recruitment and selection of personnel ; placement and monitoring; performance evaluation; training and education; evaluation measures and their importance in the selection ; practical experiences in the classroom as facilitators in the selection process . Motivation : the basic and psychological needs of human beings; motivation and conflicts; forgotten factors as motivators in the company : envy, jealousy , fear, abuse of power . Motivation Assessment.- Leadership: definition, theories and leadership development- Training and Development: definition, differentiation, stages, difficulties- Performance evaluation: definition, types, frequency, importance"

# --- 9. Avaliação section: Método / Critério / Norma de recuperação -----
$d.Content.Find.Execute("utilizar-se-á provas dissertativas com estudo de caso, para levar os alunos à maior reflexão sobre a utilização dessa ciência para o futuro engenheiro em seu trabalho cotidiano em empresas.", $false, $false, $false, $false, $false, $true, 1, $false, "Os alunos poderão participar de eventos, como feiras de ciências, ou fazer contato com escolas ou entidades estudantis cujos participantes estejam em momento de busca de trabalho. Poderão ser estudantes de curso técnico, de final de graduação ou pessoas da comunidade que participem das diversas entidades estudantis da EEL.", 2) | Out-Null

$d.Content.Find.Execute("(P1+ P2) : 2 = Média.", $false, $false, $false, $false, $false, $true, 1, $false, "Espera-se que o aluno de graduação consiga aplicar conhecimento aprendido na disciplina LOB 1031, ajudando o publico-alvo a aprimorar e ampliar sua possibilidade de empregabilidade, uma vez que o currículum vitae é exigência em praticamente todos os processos de recrutamento e seleção.Espera-se que o grupo social atendido tenha mais condições de aumentar suas possibilidades de participação em processos seletivos.", 2) | Out-Null

$d.Content.Find.Execute("Trabalho e nova avaliação para recuperação da nota necessária para aprovação.", $false, $false, $false, $false, $false, $true, 1, $false, "O aluno deverá, no último mês de aula do semestre:- fazer uma pesquisa aprofundada sobre formas e formatos de curriculum vitae- fazer contato com entidades ou organizadores de evento para oferecer e acertar detalhes da aplicação da atividade- reunir-se com o participante e desenvolver a atividade proposta: explicar do que se trata um CV, coletar as informações necessárias, elaborar conjuntamente o CV, fornecer orientações básicas para apresentação em entrevista de emprego, aplicar o questionário de satisfação", 2) | Out-Null

# --- 10. Bibliografia paragraph -> "Cada participante..." ---------------
$p19 = $d.Paragraphs(19)
$p19.Range.Text = "Cada participante atendido preencherá um formulário de satisfação, ao final da atividade."
